$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Generated Features" sheet: update the "density" row and append the new
#    generated-feature rows (variable name + description) that were added to
#    the data dictionary.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Generated Features")

# Row 3 used to be just "density" in column A with no description.
# It becomes "density_est_2018" / "Estimated density for year 2018."
$ws3.Range("B3").Value = "Estimated density for year 2018."
$ws3.Range("A3").Value = "density_est_2018"

# New rows 4-7: per-capita store count features. The variable names were
# filled down column A first, then the descriptions down column B (matching
# shared-string insertion order in the saved workbook).
$ws3.Range("A4").Value = "groc14_per_capita"
$ws3.Range("A5").Value = "superc14_per_capita"
$ws3.Range("A6").Value = "convs14_per_capita"
$ws3.Range("A7").Value = "specs14_per_capita"

$ws3.Range("B4").Value = "Count of grocery stores in county 2014 per capita."
$ws3.Range("B5").Value = "Count of supercenter stores in county 2014 per capita."
$ws3.Range("B6").Value = "Count of convenience stores in county 2014 per capita."
$ws3.Range("B7").Value = "Count of specialty food stores in county 2014 per capita."

# ---------------------------------------------------------------------------
# 2. View/selection state updates recorded by the author while reviewing the
#    workbook. Touch the other two sheets first (without leaving them as the
#    active tab) then re-activate "Generated Features" last so it stays the
#    active sheet, matching the saved workbook state.
# ---------------------------------------------------------------------------

# "Data Description" sheet: scrolled/zoomed out and left a different cell
# selected.
$ws1 = $wb.Worksheets.Item("Data Description")
$ws1.Select()
$ws1.Range("C41").Select()
$excel.ActiveWindow.Zoom = 80

# "RUCA Codes" sheet: a cell got selected on that tab too.
$ws2 = $wb.Worksheets.Item("RUCA Codes")
$ws2.Select()
$ws2.Range("B38").Select()

# "Generated Features" selection moves from B3 to A2, and the sheet remains
# the active tab in the saved workbook.
$ws3.Activate()
$ws3.Range("A2").Select()
